# Daily attendance processing - rotate "Recorded By" (column G) name lists
# so the most recently recorded-by name moves to the front of the list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) {
        continue
    }
    $parts = @($text -split ", ")
    if ($parts.Count -gt 1) {
        $lastPart = $parts[$parts.Count - 1]
        $rest = $parts[0..($parts.Count - 2)]
        $rotated = @($lastPart) + $rest
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
